$d = $word.ActiveDocument

# Paragraph 2 is the intro line ("These are the childless tags...").
# Paragraphs 3 through the end (Count) are the PUMP:* tag lines to remove.
$paras = $d.Paragraphs
$count = $paras.Count

$startPara = $paras.Item(3)
$endPara = $paras.Item($count)

$start = $startPara.Range.Start
$end = $endPara.Range.End

$r = $d.Range($start, $end)
$r.Delete()
